$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 399, shifting existing rows 399:454 down to 400:455
$ws.Rows(399).Insert()

# Populate the newly inserted row 399 with its data
$ws.Cells.Item(399, 1).Value = 11
$ws.Cells.Item(399, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(399, 3).Value = "Bíobío"
$ws.Cells.Item(399, 4).Value = 45154
$ws.Cells.Item(399, 5).Value = 8
$ws.Cells.Item(399, 6).Value = 100112009
$ws.Cells.Item(399, 7).Value = "Acelga"
$ws.Cells.Item(399, 8).Value = "Sin especificar"
$ws.Cells.Item(399, 9).Value = "Primera"
$ws.Cells.Item(399, 10).Value = 180
$ws.Cells.Item(399, 11).Value = 600
$ws.Cells.Item(399, 12).Value = 650
$ws.Cells.Item(399, 13).Value = 622
$ws.Cells.Item(399, 14).Value = "`$/atado 0,5 a 1 kilo"
$ws.Cells.Item(399, 15).Value = "Región de Ñuble"
$ws.Cells.Item(399, 16).Value = 622
$ws.Cells.Item(399, 17).Value = 1
$ws.Cells.Item(399, 18).Value = "Hortaliza"
